$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2564746666666666
$ws.Range("H2").Value = 0.7694239999999999
$ws.Range("I2").Value = 0.1818007399394835
$ws.Range("J2").Value = 0.1818007399394835
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 0.07442621253688889
$ws.Range("R2").Value = 0.6698359128319999
$ws.Range("S2").Value = 0.006234655628438898
$ws.Range("T2").Value = 0.006234655628438899

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2564746666666666
$ws.Range("H3").Value = 0.7694239999999999
$ws.Range("I3").Value = 0.1818007399394835
$ws.Range("J3").Value = 0.1818007399394835
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 1.823624817116444
$ws.Range("R3").Value = 16.412623354048
$ws.Range("S3").Value = 0.1527643600641453
$ws.Range("T3").Value = 0.1527643600641453

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2564746666666666
$ws.Range("H4").Value = 0.7694239999999999
$ws.Range("I4").Value = 0.1818007399394835
$ws.Range("J4").Value = 0.1818007399394835
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 0.2721956233262222
$ws.Range("R4").Value = 2.449760609936
$ws.Range("S4").Value = 0.02280172424689938
$ws.Range("T4").Value = 0.02280172424689937

# Row 5
$ws.Range("I5").Value = 0.7694380609030022
$ws.Range("J5").Value = 0.7694380609030022
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 0.3149952011955556
$ws.Range("R5").Value = 2.83495681076
$ws.Range("S5").Value = 0.02638702867073514
$ws.Range("T5").Value = 0.02638702867073514

# Row 6
$ws.Range("I6").Value = 0.7694380609030022
$ws.Range("J6").Value = 0.7694380609030022
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.6465469448692601
$ws.Range("T6").Value = 0.6465469448692601

# Row 7
$ws.Range("I7").Value = 0.7694380609030022
$ws.Range("J7").Value = 0.7694380609030022
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("S7").Value = 0.09650408736300695
$ws.Range("T7").Value = 0.09650408736300693

# Row 8
$ws.Range("G8").Value = 0.06878966666666667
$ws.Range("I8").Value = 0.0487611991575143
$ws.Range("J8").Value = 0.0487611991575143
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 0.01996202751022222
$ws.Range("R8").Value = 0.179658247592
$ws.Range("S8").Value = 0.001672211482076602
$ws.Range("T8").Value = 0.001672211482076602

# Row 9
$ws.Range("G9").Value = 0.06878966666666667
$ws.Range("I9").Value = 0.0487611991575143
$ws.Range("J9").Value = 0.0487611991575143
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("Q9").Value = 0.4891186522431111
$ws.Range("S9").Value = 0.04097328419970991
$ws.Range("T9").Value = 0.04097328419970991

# Row 10
$ws.Range("G10").Value = 0.06878966666666667
$ws.Range("I10").Value = 0.0487611991575143
$ws.Range("J10").Value = 0.0487611991575143
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("Q10").Value = 0.07300622100455556
$ws.Range("R10").Value = 0.657055989041
$ws.Range("S10").Value = 0.00611570347572779
$ws.Range("T10").Value = 0.006115703475727788
